$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6 (oldest reading) is dropped; rows shift up and dimension becomes A1:AH5
$ws.Rows.Item(6).Delete()

# Refresh data for the remaining rows 2-5 (new timestamps + sensor readings)
$ws.Range("A2").Value = 45107.50694444445
$ws.Range("B2").Value = 16.334
$ws.Range("C2").Value = 10.878
$ws.Range("D2").Value = 3.753
$ws.Range("E2").Value = 34.787
$ws.Range("F2").Value = 27.078
$ws.Range("G2").Value = 12.827
$ws.Range("H2").Value = 39.214
$ws.Range("I2").Value = 19.779
$ws.Range("J2").Value = 8.095000000000001
$ws.Range("K2").Value = 12.091
$ws.Range("L2").Value = 13.724
$ws.Range("M2").Value = 14.293
$ws.Range("N2").Value = 4.102
$ws.Range("O2").Value = 12.783
$ws.Range("P2").Value = 17.703
$ws.Range("Q2").Value = 11.145
$ws.Range("R2").Value = 3.152
$ws.Range("S2").Value = 2.018
$ws.Range("T2").Value = 187.2
$ws.Range("U2").Value = 35.47
$ws.Range("V2").Value = 11.799
$ws.Range("W2").Value = 23.043
$ws.Range("X2").Value = 11.665
$ws.Range("Y2").Value = 3.17
$ws.Range("Z2").Value = 20.069
$ws.Range("AA2").Value = 10.422
$ws.Range("AB2").Value = 9.335000000000001
$ws.Range("AC2").Value = 11.223
$ws.Range("AD2").Value = 14.442
$ws.Range("AE2").Value = 3.314
$ws.Range("AF2").Value = 35.145
$ws.Range("AG2").Value = 6.324
$ws.Range("AH2").Value = 14.751
$ws.Range("A3").Value = 45107.51388888889
$ws.Range("B3").Value = 2.883
$ws.Range("C3").Value = 1.486
$ws.Range("D3").Value = 1.324
$ws.Range("E3").Value = 6.124
$ws.Range("F3").Value = 4.078
$ws.Range("G3").Value = 2.244
$ws.Range("H3").Value = 13.395
$ws.Range("I3").Value = 3.49
$ws.Range("J3").Value = 1.309
$ws.Range("K3").Value = 1.614
$ws.Range("L3").Value = 2.392
$ws.Range("M3").Value = 2.379
$ws.Range("N3").Value = 0.742
$ws.Range("O3").Value = 2.256
$ws.Range("P3").Value = 3.068
$ws.Range("Q3").Value = 2.336
$ws.Range("R3").Value = 1.371
$ws.Range("S3").Value = 0.665
$ws.Range("T3").Value = 27.096
$ws.Range("U3").Value = 6.695
$ws.Range("V3").Value = 2.082
$ws.Range("W3").Value = 4.08
$ws.Range("X3").Value = 1.907
$ws.Range("Y3").Value = 0.993
$ws.Range("Z3").Value = 6.206
$ws.Range("AA3").Value = 1.839
$ws.Range("AB3").Value = 1.829
$ws.Range("AC3").Value = 2.201
$ws.Range("AD3").Value = 2.503
$ws.Range("AE3").Value = 1.246
$ws.Range("AF3").Value = 12.972
$ws.Range("AG3").Value = 0.9350000000000001
$ws.Range("AH3").Value = 2.608
$ws.Range("A4").Value = 45107.52083333334
$ws.Range("B4").Value = 16.334
$ws.Range("C4").Value = 11.823
$ws.Range("D4").Value = 1.306
$ws.Range("E4").Value = 35.473
$ws.Range("F4").Value = 28.621
$ws.Range("G4").Value = 12.829
$ws.Range("H4").Value = 45.632
$ws.Range("I4").Value = 19.779
$ws.Range("J4").Value = 8.644
$ws.Range("K4").Value = 12.743
$ws.Range("L4").Value = 14.215
$ws.Range("M4").Value = 14.935
$ws.Range("N4").Value = 4.104
$ws.Range("O4").Value = 12.783
$ws.Range("P4").Value = 18.057
$ws.Range("Q4").Value = 10.941
$ws.Range("R4").Value = 1.016
$ws.Range("S4").Value = 0.822
$ws.Range("T4").Value = 187.22
$ws.Range("U4").Value = 35.607
$ws.Range("V4").Value = 11.799
$ws.Range("W4").Value = 23.728
$ws.Range("X4").Value = 12.437
$ws.Range("Y4").Value = 2.132
$ws.Range("Z4").Value = 22.605
$ws.Range("AA4").Value = 10.422
$ws.Range("AB4").Value = 9.305999999999999
$ws.Range("AC4").Value = 10.97
$ws.Range("AD4").Value = 14.927
$ws.Range("AE4").Value = 0.766
$ws.Range("AF4").Value = 41.297
$ws.Range("AG4").Value = 6.547
$ws.Range("AH4").Value = 14.751
$ws.Range("A5").Value = 45107.52777777778
$ws.Range("B5").Value = 12.01
$ws.Range("C5").Value = 8.69
$ws.Range("D5").Value = 0.95
$ws.Range("E5").Value = 26.09
$ws.Range("F5").Value = 21.02
$ws.Range("G5").Value = 9.43
$ws.Range("H5").Value = 37.82
$ws.Range("I5").Value = 14.54
$ws.Range("J5").Value = 6.38
$ws.Range("K5").Value = 9.34
$ws.Range("L5").Value = 10.47
$ws.Range("M5").Value = 11
$ws.Range("N5").Value = 3.02
$ws.Range("O5").Value = 9.4
$ws.Range("P5").Value = 13.31
$ws.Range("Q5").Value = 8.07
$ws.Range("R5").Value = 0.76
$ws.Range("S5").Value = 0.59
$ws.Range("T5").Value = 135.74
$ws.Range("U5").Value = 26.34
$ws.Range("V5").Value = 8.68
$ws.Range("W5").Value = 17.54
$ws.Range("X5").Value = 9.16
$ws.Range("Y5").Value = 1.57
$ws.Range("Z5").Value = 18.23
$ws.Range("AA5").Value = 7.66
$ws.Range("AB5").Value = 6.86
$ws.Range("AC5").Value = 8.07
$ws.Range("AD5").Value = 10.99
$ws.Range("AE5").Value = 0.55
$ws.Range("AF5").Value = 34.58
$ws.Range("AG5").Value = 4.81
$ws.Range("AH5").Value = 10.85

# A handful of columns grew from 7 to 8 characters wide (custom accuracy formatting)
$ws.Columns.Item(2).ColumnWidth = 7.17
$ws.Columns.Item(3).ColumnWidth = 7.17
$ws.Columns.Item(5).ColumnWidth = 7.17
$ws.Columns.Item(6).ColumnWidth = 7.17
$ws.Columns.Item(7).ColumnWidth = 7.17
$ws.Columns.Item(8).ColumnWidth = 7.17
$ws.Columns.Item(9).ColumnWidth = 7.17
$ws.Columns.Item(11).ColumnWidth = 7.17
$ws.Columns.Item(12).ColumnWidth = 7.17
$ws.Columns.Item(13).ColumnWidth = 7.17
$ws.Columns.Item(15).ColumnWidth = 7.17
$ws.Columns.Item(16).ColumnWidth = 7.17
$ws.Columns.Item(17).ColumnWidth = 7.17
$ws.Columns.Item(21).ColumnWidth = 7.17
$ws.Columns.Item(22).ColumnWidth = 7.17
$ws.Columns.Item(23).ColumnWidth = 7.17
$ws.Columns.Item(24).ColumnWidth = 7.17
$ws.Columns.Item(26).ColumnWidth = 7.17
$ws.Columns.Item(27).ColumnWidth = 7.17
$ws.Columns.Item(29).ColumnWidth = 7.17
$ws.Columns.Item(30).ColumnWidth = 7.17
$ws.Columns.Item(32).ColumnWidth = 7.17
$ws.Columns.Item(34).ColumnWidth = 7.17
